$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: "Applicant" -> "Case title", "representative" -> "Representative"
$ws.Range("B1").Value = "Case title"
$ws.Range("C1").Value = "Representative"

# Update data rows: "Applicant A" -> "Case A", "Applicant B" -> "Case B"
$ws.Range("B2").Value = "Case A"
$ws.Range("B3").Value = "Case B"

# Update the selected cell to match the saved view state
$ws.Range("D18").Select()
